# Update the pending-documents tracking sheet:
#  - row 3 is replaced by a new document entry for order P-22/074-S00
#    (doc 22-074-PLG-0005), pushing the former row-3 record down to row 5
#  - row 6 (order P-23/028-S00) is updated with a new document/revision
#    (doc 23-028-PRC-0012) replacing its previous document/revision data
#  - row 7 keeps order P-23/036-S00 (shifted down from the former row 6)
#  - the former row 7 (order P-23/048-S00) is dropped
#  - a new row 8 is appended for order P-23/044-S05 (doc 23-044-S05-PLG-0005-B18)
# Net effect: sheet grows from A1:P7 to A1:P8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (new content: P-22/074-S00 / 22-074-PLG-0005)
$ws.Range("A3").Value = "P-22/074-S00"
$ws.Range("B3").Value = "22-12-2022"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "10-07-2023"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "104001091"
$ws.Range("F3").Value = "Caudal"
$ws.Range("G3").Value = "V-1040010910-0002"
$ws.Range("H3").Value = "22-074-PLG-0005"
$ws.Range("I3").Value = "OVERALL DRAWING WITH WEIGHT FOR FE & FO (ARZANAH)"
$ws.Range("J3").Value = "Planos"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "2"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "12-07-2024"
$ws.Range("P3").Value = "05-01-2024 Enviado Rev. 0 // 22-02-2024 Com. Mayores Rev. 1 // 22-03-2024 Enviado Rev. 1 // 22-05-2024 Com. Menores Rev. 1 // 29-05-2024 Enviado Rev. 2 // 12-07-2024 Com. Menores Rev. 2"

# Row 5 (former row-3 content: P-22/075-S00 / 22-075-PLN-0001)
$ws.Range("A5").Value = "P-22/075-S00"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "09-12-2022"
$ws.Range("C5").Value = "27-06-2023"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "104001071"
$ws.Range("F5").Value = "Nivel"
$ws.Range("G5").Value = "V-1040010710-0003"
$ws.Range("H5").Value = "22-075-PLN-0001"
$ws.Range("I5").Value = "QUALITY CONTROL PLAN 22-075"
$ws.Range("J5").Value = "PPI"
$ws.Range("K5").Value = "Sí"
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "5"
$ws.Range("N5").Value = "24-05-2024"
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = "31-10-2023 Comentado Rev. 2 // 07-11-2023 Enviado Rev. 3 // 21-12-2023 Com. Menores Rev. 3 // 05-01-2024 Enviado Rev. 4 // 07-02-2024 Com. Menores Rev. 5 // 14-02-2024 Enviado Rev. 5 // 24-05-2024 Com. Menores Rev. 5"

# Row 6 (order P-23/028-S00, new document/revision 23-028-PRC-0012)
$ws.Range("A6").Value = "P-23/028-S00"
$ws.Range("B6").Value = "28-02-2023"
$ws.Range("C6").Value = "26-10-2023"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103701061"
$ws.Range("F6").Value = "Temperatura"
$ws.Range("G6").Value = "3998_18-1037010610-00014"
$ws.Range("H6").Value = "23-028-PRC-0012"
$ws.Range("I6").Value = "PACKING & TRANSPORTATION PROCEDURE"
$ws.Range("J6").Value = "Packing"
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "1"
$ws.Range("N6").Value = "17-07-2024"
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = "11-12-2023 Com. Menores Rev. 0 // 15-02-2024 Enviado Rev. 0 // 11-06-2024 Com. Menores Rev. 0 // 11-06-2024 Enviado Rev. 1 // 17-07-2024 Com. Menores Rev. 1"

# Row 7 (former row-6 content: P-23/036-S00 / 23-036-DOS-0002)
$ws.Range("A7").Value = "P-23/036-S00"
$ws.Range("B7").Value = "28-03-2023"
$ws.Range("C7").Value = "23-11-2023"
$ws.Range("D7").Value = "RFQ 12-99-52-1807 _REV.A"
$ws.Range("F7").Value = "Caudal"
$ws.Range("G7").Value = "8005710911-V-0011"
$ws.Range("H7").Value = "23-036-DOS-0002"
$ws.Range("I7").Value = "FINAL QUALITY DOSSIER"
$ws.Range("N7").Value = "14-06-2024"
$ws.Range("O7").Value = "Este pedido esta terminado. Mientras no reclamen no vamos ha enviar nada. Entra a fecha 14/06/2024 Aceptado con Com.Menores"
$ws.Range("P7").Value = "24-07-2023 Aprobado Rev. 0 // 14-06-2024 Com. Menores Rev. 0"

# Row 8 (brand new row: P-23/044-S05 / 23-044-S05-PLG-0005-B18)
$ws.Range("A8").Value = "P-23/044-S05"
$ws.Range("B8").Value = "31-05-2024"
$ws.Range("C8").Value = "22-07-2024"
$ws.Range("D8").Value = "1037010910-05"
$ws.Range("E8").Value = "TÉCNICAS REUNIDAS"
$ws.Range("F8").Value = "Caudal"
$ws.Range("G8").Value = "3998_18-1037010910-00051"
$ws.Range("H8").Value = "23-044-S05-PLG-0005-B18"
$ws.Range("I8").Value = "NFXP3 - BARZAN ISBL - CALCULATIONS AND OVERALL DRAWINGS FOR RESTRICTION ORIFICE"
$ws.Range("J8").Value = "Cálculo y plano"
$ws.Range("K8").Value = "Sí"
$ws.Range("L8").Value = "Com. Menores"
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "1"
$ws.Range("N8").Value = "17-07-2024"
$ws.Range("O8").Value = ""
$ws.Range("P8").Value = "24-05-2024 Enviado Rev. 0 // 25-06-2024 Com. Menores Rev. 0 // 26-06-2024 Enviado Rev. 1 // 17-07-2024 Com. Menores Rev. 1"
